# Swap the contents of columns C (codeforiati:group-name) and
# D (codeforiati:group-code), including the header row, for every
# row of data in the sheet. This matches the reordering performed
# upstream in the codeforIATI/codelists source data, which now lists
# the group-code column before the group-name column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp

$rangeC = $ws.Range("C1:C$lastRow")
$rangeD = $ws.Range("D1:D$lastRow")

$valuesC = $rangeC.Value2
$valuesD = $rangeD.Value2

$rangeC.Value2 = $valuesD
$rangeD.Value2 = $valuesC
